$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("c_stocks")

# Append " Forest" to the land-use label in column E for the AGB/BGB blocks
# (rows 3-7, 9-13, 15-19, 21-25, 27-31) and the degraded-forest DG_ratio block
# (rows 39-43).
$rows = @(3,4,5,6,7,9,10,11,12,13,15,16,17,18,19,21,22,23,24,25,27,28,29,30,31,39,40,41,42,43)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value() + " Forest"
}

# Update the sheet view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E45").Select()
